$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.026.63"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.361.40"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.90"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.43"
$ws.Range("E7").Value = "  +2.18%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +11.24%  "
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.22"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.38"
$ws.Range("E12").Value = "  +12.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.33"
$ws.Range("E13").Value = "  +10.25%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "2.712.96"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.69"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.902"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "2.367.13"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "43.921.22"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  +5.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.10"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.62"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +26.43%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.66"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.79"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  -3.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.79"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.39"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0769"
$ws.Range("E34").Value = "  +7.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +4.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.76"
$ws.Range("E37").Value = "  -5.73%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.33"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0279"
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("E41").Value = "  +15.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.207"
$ws.Range("E42").Value = "  +14.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.14"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.09"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.75"
$ws.Range("E46").Value = "  +6.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.52"
$ws.Range("E47").Value = "  +11.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.51"
$ws.Range("E48").Value = "  +8.40%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.51"
$ws.Range("E51").Value = "  +2.47%  "
